$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "usuario conexiones totales": add a new "Cargo" column (E) of data
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("usuario conexiones totales")

$ws2.Range("E2").Value = "prueba"
$ws2.Range("E3").Value = "CEO"
# Row 4 keeps the same "blank" text value the rest of that row already uses
# (column D on that row is an empty string) - force a literal empty string
# value via a formula so the cell keeps its Text type instead of being
# cleared outright.
$ws2.Range("E4").Formula = '=""'

# ---------------------------------------------------------------------------
# Sheet "Inscritos a evento": add "Cargo" (G) and "Hora de inscripcion" (H)
# columns of data for every attendee row
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Inscritos a evento")

# 2020-12-16 12:47:00, expressed as the Excel date serial number (days since
# 1899-12-30) so the stored value matches exactly.
$horaInscripcion = 44181.5326388889

$cargos = @{
    2 = "CEO"
    3 = "Cluster Manager"
    4 = ""
    5 = ""
    6 = "Producer"
    7 = "prueba"
    8 = "Lead Developer"
}

foreach ($row in 2..8) {
    $cargo = $cargos[$row]
    if ($cargo -eq "") {
        $ws4.Range("G$row").Formula = '=""'
    } else {
        $ws4.Range("G$row").Value = $cargo
    }

    $ws4.Range("H$row").Value = $horaInscripcion
    $ws4.Range("H$row").NumberFormat = "yyyy-mm-dd h:mm:ss"
}
